# Rerun of the jig calibration matrices: the "x" / "y" curvature-vs-signal
# blocks (rows 4-8 / 10-14) on all three AA1/AA2/AA3 sheets get replaced
# with the new (tiny, negative) calibration readings, and the view state
# (selection / active sheet) moves on to reflect where the author ended up
# looking after rerunning the calc.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "AA1"
$ws2 = $wb.Worksheets.Item(2)   # "AA2"
$ws3 = $wb.Worksheets.Item(3)   # "AA3"

# New calibration values (exact doubles matching the target workbook).
$v1 = -0.0000903954802259886993840
$v2 = -0.0001412429378531073495637
$v3 = -0.0001807909604519773987680
$v4 = -0.0002259887005649717484600
$v5 = -0.0002824858757062146991275

# ---- AA1 --------------------------------------------------------------
# Column A (rows 4-8) goes from 0.8/1.25/1.6/2/2.5 to the new values;
# column B (rows 4-8) is already 0 and stays that way.
$ws1.Range("A4").Value = $v1
$ws1.Range("A5").Value = $v2
$ws1.Range("A6").Value = $v3
$ws1.Range("A7").Value = $v4
$ws1.Range("A8").Value = $v5

# Column B (rows 10-14) goes from 0.8/1.25/1.6/2/2.5 to the new values;
# column A (rows 10-14) is already 0 and stays that way.
$ws1.Range("B10").Value = $v1
$ws1.Range("B11").Value = $v2
$ws1.Range("B12").Value = $v3
$ws1.Range("B13").Value = $v4
$ws1.Range("B14").Value = $v5

# ---- AA2 --------------------------------------------------------------
$ws2.Range("A4").Value = $v1
$ws2.Range("A5").Value = $v2
$ws2.Range("A6").Value = $v3
$ws2.Range("A7").Value = $v4
$ws2.Range("A8").Value = $v5

$ws2.Range("B4").Value = 0
$ws2.Range("B5").Value = 0
$ws2.Range("B6").Value = 0
$ws2.Range("B7").Value = 0
$ws2.Range("B8").Value = 0

$ws2.Range("A10").Value = 0
$ws2.Range("A11").Value = 0
$ws2.Range("A12").Value = 0
$ws2.Range("A13").Value = 0
$ws2.Range("A14").Value = 0

$ws2.Range("B10").Value = $v1
$ws2.Range("B11").Value = $v2
$ws2.Range("B12").Value = $v3
$ws2.Range("B13").Value = $v4
$ws2.Range("B14").Value = $v5

# ---- AA3 --------------------------------------------------------------
$ws3.Range("A4").Value = $v1
$ws3.Range("A5").Value = $v2
$ws3.Range("A6").Value = $v3
$ws3.Range("A7").Value = $v4
$ws3.Range("A8").Value = $v5

$ws3.Range("B4").Value = 0
$ws3.Range("B5").Value = 0
$ws3.Range("B6").Value = 0
$ws3.Range("B7").Value = 0
$ws3.Range("B8").Value = 0

$ws3.Range("A10").Value = 0
$ws3.Range("A11").Value = 0
$ws3.Range("A12").Value = 0
$ws3.Range("A13").Value = 0
$ws3.Range("A14").Value = 0

$ws3.Range("B10").Value = $v1
$ws3.Range("B11").Value = $v2
$ws3.Range("B12").Value = $v3
$ws3.Range("B13").Value = $v4
$ws3.Range("B14").Value = $v5

# ---- View state ---------------------------------------------------------
# AA1 keeps its A3:B14 selection but is no longer the active tab.
$null = $ws1.Activate()
$null = $ws1.Range("A3:B14").Select()

# AA2's selection moves from A24:C29 up to A3:B14 (and its old A7
# top-left scroll position resets along with it).
$null = $ws2.Activate()
$null = $ws2.Range("A3:B14").Select()

# AA3 becomes the active / selected tab, with the selection moved to C5.
$null = $ws3.Activate()
$null = $ws3.Range("C5").Select()
